$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Property1" to "DataNode"
$ws.Name = "DataNode"

# Move selection to C24 (reflects cursor position at save time)
$ws.Range("C24").Select()
